$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Produto")
$v = $ws.Range("A1").Value
Write-Host "VALUE: $v"
Write-Host "TYPE: $($v.GetType())"
